# revert test file changes and fix build
#
# The "10am" / "1:30pm" time values in column C had trailing padding
# whitespace baked into the shared strings; strip it back to the clean
# values. Also restore the saved cell-selection to C2 (it had drifted to
# C6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "10am"
$ws.Range("C3").Value = "1:30pm"

$ws.Range("C2").Select()
